$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.667686462402344
$ws.Range("B1").Value = 2.444854974746704
$ws.Range("C1").Value = 2.075833320617676
$ws.Range("D1").Value = 1.5724858045578
$ws.Range("E1").Value = 1.241879940032959
